$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New downtime-log rows (Line:9 / Line:8 events) appended below the header,
# feeding the highcharts visualization added for this line file.
$rows = @(
    @{ Row=2; A="Line:9 Stage:1"; B="01/09/2024"; C="pri cl LA"; D=45300.40887037037; E=45300.40910185185; F=0.33; G="44bb2153-845f-4543-8e0b-e127667e7e30"; H=0.0002314814814814815 }
    @{ Row=3; A="Line:9 Stage:1"; B="01/09/2024"; C="Pri pH flows"; D=45300.49280315972; E=45300.49893741898; F=8.83; G="44bb2153-845f-4543-8e0b-e127667e7e30"; H=0.006134259259259259 }
    @{ Row=4; A="Line:9 Stage:1"; B="01/09/2024"; C="Pri pH flows"; D=45300.50241207176; E=45300.50287503472; F=0.67; G="44bb2153-845f-4543-8e0b-e127667e7e30"; H=0.000462962962962963 }
    @{ Row=5; A="Line:9 Stage:1"; B="01/09/2024"; C="Pri pH flows"; D=45300.50356979167; E=45300.50368553241; F=0.17; G="44bb2153-845f-4543-8e0b-e127667e7e30"; H=0.0001157407407407407 }
    @{ Row=6; A="Line:9 Stage:1"; B="01/09/2024"; C="Pri pH flows"; D=45300.50438012731; E=45300.50449586806; F=0.17; G="44bb2153-845f-4543-8e0b-e127667e7e30"; H=0.0001157407407407407 }
    @{ Row=7; A="Line:9 Stage:1"; B="01/09/2024"; C="Pri pH flows"; D=45300.50519054398; E=45300.50715813658; F=2.83; G="44bb2153-845f-4543-8e0b-e127667e7e30"; H=0.001967592592592592 }
    @{ Row=8; A="Line:8 Stage:1"; B="01/09/2024"; C="pri cl LA"; D=45300.42056299769; E=45300.42067873842; F=0.17; G="44bb2153-845f-4543-8e0b-e127667e7e30"; H=0.0001157407407407407 }
    @{ Row=9; A="Line:8 Stage:1"; B="01/09/2024"; C="pri pH HA"; D=45300.65476473379; E=45300.65488047454; F=0.17; G="44bb2153-845f-4543-8e0b-e127667e7e30"; H=0.0001157407407407407 }

)

# --- Pass 1: write the "st"/"nd" timestamp columns (D, E) first so the
#     datetime number-format style is created once and reused for every
#     row (matches a single shared cellXfs entry instead of one per row).
foreach ($r in $rows) {
    if ($r.Row -eq 2) {
        $ws.Range("D2").NumberFormat = "yyyy-mm-dd h:mm:ss"
    }
    $ws.Range("D$($r.Row)").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r.Row, 4).Value = $r.D

    $ws.Range("E$($r.Row)").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

# --- Pass 2: write the remaining text/numeric columns (A, B, C, F, G, H).
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A

    # "tdate" (B) is a plain text column, but values like "01/09/2024" look
    # like dates, so force text formatting before typing, then drop back to
    # the default style once the literal text is stored.
    $ws.Range("B$($r.Row)").NumberFormat = "@"
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Range("B$($r.Row)").Style = "Normal"

    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
}

Write-Output "done"
